$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 absorbs the data that used to live in row 3 (the duplicate
# FAPs/Crp/Olr1/FAPs self-edge), and the derived-specificity columns are
# recalculated now that only a single edge remains. Row 3 is then removed.
$ws.Range("D2").Value = "FAPs"

$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.006582333333333333
$ws.Range("N2").Value = 0.019747
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.0005043844563333333
$ws.Range("R2").Value = 0.004539460107
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Remove the now-redundant row 3 entirely.
$ws.Rows("3").Delete()
